# Adds the new "MES SETUP" panel translation rows (66-71) to the
# "Translation" sheet, matching the TouchGFX texts.xlsx commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$rows = @(
    @{ Row = 66; B = "SingleUseId66"; C = "Default"; D = "Left";   E = "LTR"; F = "HF INPUT" },
    @{ Row = 67; B = "SingleUseId67"; C = "Large";   D = "Left";   E = "LTR"; F = "ON" },
    @{ Row = 68; B = "SingleUseId68"; C = "Large";   D = "Left";   E = "LTR"; F = "OFF" },
    @{ Row = 69; B = "SingleUseId69"; C = "Default"; D = "Left";   E = "LTR"; F = "GATE" },
    @{ Row = 70; B = "SingleUseId70"; C = "Default"; D = "Center"; E = "LTR"; F = "<value> ms" },
    @{ Row = 71; B = "SingleUseId71"; C = "Default"; D = "Left";   E = "LTR"; F = "MES SETUP" }
)

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
}
